$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the exogenous agent parameter example values on row 3
# (man_num, woman_num, selector_num, candidate_num)
$ws.Range("D3").Value = 15
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 10

# Reflect the new active selection / zoom used when the sample values were set
$ws.Range("D4").Select()
$ws.Application.ActiveWindow.Zoom = 235
